$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

# Update column A values (rows 2-12) from "2" to "24" (kept as text)
for ($r = 2; $r -le 12; $r++) {
    $ws.Cells.Item($r, 1).Value = "'24"
}

# Delete the entire column M (was "Event" header + "nan"/empty data)
$ws.Range("M1:M12").EntireColumn.Delete()
